$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Vi V" wallet column is being inserted into the customer list header
# (row 2) just before "Trang thai", pushing "Trang thai" and "Ngay tao" one
# column to the right (G -> H -> I).
$ngayTao  = $ws.Range("H2").Value2
$trangThai = $ws.Range("G2").Value2

$ws.Range("I2").Value2 = $ngayTao
$ws.Range("H2").Value2 = $trangThai
$ws.Range("G2").Value2 = "Ví V"

$null = $ws.Range("G6").Select()
